$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("H2").Value = 0.1069567095850491
$ws.Range("B3").Value = 0.09043386850983193
$ws.Range("H3").Value = 0.1973905780948811
$ws.Range("B4").Value = 0.08657203724954163
$ws.Range("H4").Value = 0.1935287468345908
$ws.Range("B5").Value = 0.07076940770021801
$ws.Range("H5").Value = 0.1777261172852672
$ws.Range("B6").Value = 0.08794861539697796
$ws.Range("H6").Value = 0.1949053249820271
$ws.Range("B7").Value = 0.01160641164542359
$ws.Range("C7").Value = 0.002110497307708295
$ws.Range("D7").Value = 1.469620788740304
$ws.Range("E7").Value = 0.02056724542427917
$ws.Range("F7").Value = 0.007438496073986227
$ws.Range("G7").Value = 0.0157743272168608
$ws.Range("H7").Value = 0.1185631212304727
$ws.Range("B8").Value = 0.009536076865655638
$ws.Range("C8").Value = 0.001926331599393779
$ws.Range("D8").Value = 1.085580795335966
$ws.Range("E8").Value = 0.03023052421355057
$ws.Range("F8").Value = 0.005754940319859566
$ws.Range("G8").Value = 0.01331721341145166
$ws.Range("H8").Value = 0.1164927864507048
$ws.Range("B9").Value = 0.00727739748016039
$ws.Range("C9").Value = 0.001764429781329991
$ws.Range("D9").Value = 0.7613382333757169
$ws.Range("E9").Value = 0.01433506537452428
$ws.Range("F9").Value = 0.003815499890623111
$ws.Range("G9").Value = 0.01073929506969786
$ws.Range("H9").Value = 0.1142341070652095
$ws.Range("B10").Value = 0.008847981047940032
$ws.Range("C10:G10").ClearContents()
$ws.Range("H10").Value = 0.1158046906329892
$ws.Range("B11").Value = 0.0270930316302893
$ws.Range("H11").Value = 0.1340497412153384
$ws.Range("B12").Value = 0.04845712219404096
$ws.Range("H12").Value = 0.1554138317790901
$ws.Range("B13").Value = 0.05989195624476158
$ws.Range("H13").Value = 0.1668486658298107
$ws.Range("B14").Value = 0.06825378053622953
$ws.Range("H14").Value = 0.1752104901212787
$ws.Range("B15").Value = 0.0740785428087717
$ws.Range("H15").Value = 0.1810352523938208
$ws.Range("B16").Value = 0.07898793123304286
$ws.Range("H16").Value = 0.185944640818092
$ws.Range("B17").Value = 0.07957936531373674
$ws.Range("H17").Value = 0.1865360748987859
$ws.Range("B18").Value = -0.1069567095850491
$ws.Range("C18").Value = 0.01087930151352742
$ws.Range("D18").Value = -18.6302705007023
$ws.Range("E18").Value = 0.03755329832689248
$ws.Range("F18").Value = -0.1283871727892488
$ws.Range("G18").Value = -0.08552624638084962
$ws.Range("H18").Value = 0
$ws.Range("B19").Value = 0.08246171007332959
$ws.Range("H19").Value = 0.1894184196583787
$ws.Range("B20").Value = 0.08443691597725367
$ws.Range("H20").Value = 0.1913936255623028
$ws.Range("B21").Value = 0.08744852886512905
$ws.Range("C21").Value = 0.00827686590926469
$ws.Range("D21").Value = 20.10065263639976
$ws.Range("E21").Value = 0.05494994659393142
$ws.Range("F21").Value = 0.0711753566427236
$ws.Range("G21").Value = 0.1037217010875345
$ws.Range("H21").Value = 0.1944052384501782
$ws.Range("B22").Value = 0.09399495848337078
$ws.Range("H22").Value = 0.2009516680684199
$ws.Range("B23").Value = 0.0956437197288324
$ws.Range("H23").Value = 0.2026004293138816
$ws.Range("B24").Value = 0.09933894505425432
$ws.Range("C24").Value = 0.008991469368418069
$ws.Range("D24").Value = 21.88141561751501
$ws.Range("E24").Value = 0.05500252448082291
$ws.Range("F24").Value = 0.08163086045605795
$ws.Range("G24").Value = 0.1170470296524511
$ws.Range("H24").Value = 0.2062956546393035
$ws.Range("B25").Value = 0.1002298850499006
$ws.Range("C25").Value = 0.008735023889623637
$ws.Range("D25").Value = 22.06344703063166
$ws.Range("E25").Value = 0.05725570524262481
$ws.Range("F25").Value = 0.08304093042708108
$ws.Range("G25").Value = 0.1174188396727207
$ws.Range("H25").Value = 0.2071865946349497
$ws.Range("B26").Value = 0.09919511808624296
$ws.Range("C26").Value = 0.009216307508984183
$ws.Range("D26").Value = 21.33043148661273
$ws.Range("E26").Value = 0.0663620139429116
$ws.Range("F26").Value = 0.0810468976343559
$ws.Range("G26").Value = 0.1173433385381295
$ws.Range("H26").Value = 0.2061518276712921
$ws.Range("B27").Value = 0.1022434263271316
$ws.Range("C27").Value = 0.009302710074219846
$ws.Range("D27").Value = 20.04756509900071
$ws.Range("E27").Value = 0.07615426964768562
$ws.Range("F27").Value = 0.0839287846934087
$ws.Range("G27").Value = 0.1205580679608547
$ws.Range("H27").Value = 0.2092001359121807
$ws.Range("B28").Value = 0.09944184158901578
$ws.Range("C28").Value = 0.009319826563435848
$ws.Range("D28").Value = 18.55880300867425
$ws.Range("E28").Value = 0.09409764082176228
$ws.Range("F28").Value = 0.08106222313164926
$ws.Range("G28").Value = 0.1178214600463839
$ws.Range("H28").Value = 0.2063985511740649
$ws.Range("B29").Value = 0.01755680487140073
$ws.Range("C29").Value = 0.005029041937002473
$ws.Range("D29").Value = 1.490034156925886
$ws.Range("E29").Value = 0.01626332096562937
$ws.Range("F29").Value = 0.007356858771473115
$ws.Range("G29").Value = 0.02775675097132797
$ws.Range("H29").Value = 0.1245135144564499

Write-Output "edit applied"
